# Refresh the crypto price / 1h-volume columns (D:E) with the latest
# scrape values, as produced by the scheduled GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A couple of the new prices are numeric-looking strings whose exact
# text (trailing zero / significant digit) would be lost if Excel
# auto-converted them to a real number (e.g. "2.10" -> 2.1). Mark just
# those cells as Text first so the literal digits are preserved, then
# drop the formatting override again so the cell is left unstyled,
# same as every other data cell in the sheet.
$protectedCells = @('D37', 'D42')
foreach ($cellRef in $protectedCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '43.964.13'
$ws.Range('E2').Value = '  +2.36%  '
$ws.Range('D3').Value = '2.218.95'
$ws.Range('E3').Value = '  +0.21%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '262.78'
$ws.Range('E5').Value = '  +2.14%  '
$ws.Range('D6').Value = '87.23'
$ws.Range('E6').Value = '  +13.63%  '
$ws.Range('D7').Value = '0.624'
$ws.Range('E7').Value = '  +2.35%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.608'
$ws.Range('E9').Value = '  +1.99%  '
$ws.Range('D10').Value = '45.79'
$ws.Range('E10').Value = '  +8.39%  '
$ws.Range('D11').Value = '0.0923'
$ws.Range('E11').Value = '  +2.00%  '
$ws.Range('E12').Value = '  +8.53%  '
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').Value = '2.549.50'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '14.66'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').Value = '2.209.50'
$ws.Range('E16').Value = '  -0.12%  '
$ws.Range('D17').Value = '0.788'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '43.937.35'
$ws.Range('E18').Value = '  +2.53%  '
$ws.Range('E19').Value = '  +1.82%  '
$ws.Range('D20').Value = '5.99'
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = '70.05'
$ws.Range('E21').Value = '  -1.79%  '
$ws.Range('D22').Value = '2.38'
$ws.Range('E22').Value = '  +7.57%  '
$ws.Range('D23').Value = '232.21'
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').Value = '9.01'
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('E26').Value = '  -0.47%  '
$ws.Range('D27').Value = '3.54'
$ws.Range('E27').Value = '  +5.88%  '
$ws.Range('D28').Value = '39.85'
$ws.Range('E28').Value = '  -6.42%  '
$ws.Range('E29').Value = '  +2.97%  '
$ws.Range('D30').Value = '2.22'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('D31').Value = '174.88'
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').Value = '20.59'
$ws.Range('E32').Value = '  +1.36%  '
$ws.Range('D33').Value = '0.0885'
$ws.Range('E33').Value = '  +1.33%  '
$ws.Range('D34').Value = '5.44'
$ws.Range('E34').Value = '  +4.04%  '
$ws.Range('D36').Value = '0.112'
$ws.Range('E36').Value = '  +4.68%  '
$ws.Range('D37').Value = '0.0360'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').Value = '4.47'
$ws.Range('E38').Value = '  +3.21%  '
$ws.Range('D39').Value = '3.31'
$ws.Range('E39').Value = '  +15.81%  '
$ws.Range('D40').Value = '12.47'
$ws.Range('E40').Value = '  -3.00%  '
$ws.Range('D41').Value = '65.24'
$ws.Range('E41').Value = '  +8.61%  '
$ws.Range('D42').Value = '2.10'
$ws.Range('D43').Value = '5.57'
$ws.Range('D44').Value = '0.203'
$ws.Range('E44').Value = '  +1.90%  '
$ws.Range('D45').Value = '101.47'
$ws.Range('E45').Value = '  -1.47%  '
$ws.Range('D46').Value = '8.37'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').Value = '0.0984'
$ws.Range('E47').Value = '  +0.67%  '
$ws.Range('E48').Value = '  +1.21%  '
$ws.Range('E49').Value = '  +4.54%  '
$ws.Range('D50').Value = '0.446'
$ws.Range('E50').Value = '  -2.53%  '
$ws.Range('D51').Value = '1.53'
$ws.Range('E51').Value = '  +7.39%  '

foreach ($cellRef in $protectedCells) {
    $ws.Range($cellRef).ClearFormats()
}
